$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TicketLog4thS")

$ws.Cells.Item(94,1).Value = "DHTD-104"
$ws.Cells.Item(94,2).Value = "Analysis & Design"
$ws.Cells.Item(94,3).Value = "Construction"
$ws.Cells.Item(94,4).Value = "Create design for car enemy"
$ws.Cells.Item(94,6).Value = 3
$ws.Cells.Item(94,7).Value = 3.5
$ws.Cells.Item(94,8).Value = "Nicolas"
$ws.Cells.Item(94,9).Value = 3.5
$ws.Range("D94").WrapText = $true
$ws.Range("H94").WrapText = $true
$ws.Range("I94").WrapText = $true

$ws.Cells.Item(95,1).Value = "DHTD-107"
$ws.Cells.Item(95,2).Value = "Project Management"
$ws.Cells.Item(95,3).Value = "Transition"
$ws.Cells.Item(95,4).Value = "Update Time Sheet (Sprint 11-16)"
$ws.Cells.Item(95,6).Value = 1
$ws.Cells.Item(95,7).Value = 1.25
$ws.Cells.Item(95,8).Value = "Nicolas"
$ws.Cells.Item(95,9).Value = 1.25
$ws.Range("D95").WrapText = $true
$ws.Range("H95").WrapText = $true
$ws.Range("I95").WrapText = $true

$ws.Cells.Item(96,1).Value = "DHTD-108"
$ws.Cells.Item(96,2).Value = "Project Management"
$ws.Cells.Item(96,3).Value = "Construction"
$ws.Cells.Item(96,4).Value = "Blog Entry Week 17"
$ws.Cells.Item(96,6).Value = 2.5
$ws.Cells.Item(96,7).Value = 3
$ws.Cells.Item(96,8).Value = "Fabian"
$ws.Cells.Item(96,9).Value = 0.5
$ws.Range("D96").WrapText = $true
$ws.Range("H96").WrapText = $true
$ws.Range("I96").WrapText = $true

$ws.Cells.Item(97,8).Value = "Luca"
$ws.Cells.Item(97,9).Value = 0.5
$ws.Range("H97").WrapText = $true
$ws.Range("I97").WrapText = $true

$ws.Cells.Item(98,8).Value = "Nicolas"
$ws.Cells.Item(98,9).Value = 2
$ws.Range("H98").WrapText = $true
$ws.Range("I98").WrapText = $true

$ws.Cells.Item(99,1).Value = "DHTD-109"
$ws.Cells.Item(99,2).Value = "Project Management"
$ws.Cells.Item(99,3).Value = "Construction"
$ws.Cells.Item(99,4).Value = "Write comments to other groups (Week 17)"
$ws.Cells.Item(99,6).Value = 1
$ws.Cells.Item(99,7).Value = 1
$ws.Cells.Item(99,8).Value = "Luca"
$ws.Cells.Item(99,9).Value = 1
$ws.Range("D99").WrapText = $true
$ws.Range("H99").WrapText = $true
$ws.Range("I99").WrapText = $true

$ws.Cells.Item(100,1).Value = "DHTD-110"
$ws.Cells.Item(100,2).Value = "Implementation"
$ws.Cells.Item(100,3).Value = "Construction"
$ws.Cells.Item(100,4).Value = "Create levels for towers"
$ws.Cells.Item(100,5).Value = "upgrade towers"
$ws.Cells.Item(100,6).Value = 2
$ws.Cells.Item(100,7).Value = 1.25
$ws.Cells.Item(100,8).Value = "Luca"
$ws.Cells.Item(100,9).Value = 1.25
$ws.Range("D100").WrapText = $true
$ws.Range("H100").WrapText = $true
$ws.Range("I100").WrapText = $true

$ws.Cells.Item(101,1).Value = "DHTD-111"
$ws.Cells.Item(101,2).Value = "Implementation"
$ws.Cells.Item(101,3).Value = "Construction"
$ws.Cells.Item(101,4).Value = "Lasertower: Check if it is posssible to make multiple damage"
$ws.Cells.Item(101,6).Value = 3
$ws.Cells.Item(101,7).Value = 1.5
$ws.Cells.Item(101,8).Value = "Fabian"
$ws.Cells.Item(101,9).Value = 1
$ws.Range("D101").WrapText = $true
$ws.Range("H101").WrapText = $true
$ws.Range("I101").WrapText = $true

$ws.Cells.Item(102,8).Value = "Luca"
$ws.Cells.Item(102,9).Value = 0.25
$ws.Range("H102").WrapText = $true
$ws.Range("I102").WrapText = $true

$ws.Cells.Item(103,8).Value = "Nicolas"
$ws.Cells.Item(103,9).Value = 0.25
$ws.Range("H103").WrapText = $true
$ws.Range("I103").WrapText = $true

$ws.Cells.Item(104,1).Value = "DHTD-112"
$ws.Cells.Item(104,2).Value = "Testing"
$ws.Cells.Item(104,3).Value = "Construction"
$ws.Cells.Item(104,4).Value = "Refactor code snippets according to Metrics"
$ws.Cells.Item(104,6).Value = 6
$ws.Cells.Item(104,7).Value = 5.25
$ws.Cells.Item(104,8).Value = "Nicolas"
$ws.Cells.Item(104,9).Value = 5.25
$ws.Range("D104").WrapText = $true
$ws.Range("H104").WrapText = $true
$ws.Range("I104").WrapText = $true

$ws.Cells.Item(105,1).Value = "DHTD-113"
$ws.Cells.Item(105,2).Value = "Testing"
$ws.Cells.Item(105,3).Value = "Construction"
$ws.Cells.Item(105,4).Value = "Improve Test Automation"
$ws.Cells.Item(105,6).Value = 8
$ws.Cells.Item(105,7).Value = 6.25
$ws.Cells.Item(105,8).Value = "Fabian"
$ws.Cells.Item(105,9).Value = 0.5
$ws.Range("D105").WrapText = $true
$ws.Range("H105").WrapText = $true
$ws.Range("I105").WrapText = $true

$ws.Cells.Item(106,8).Value = "Luca"
$ws.Cells.Item(106,9).Value = 0.5
$ws.Range("H106").WrapText = $true
$ws.Range("I106").WrapText = $true

$ws.Cells.Item(107,8).Value = "Nicolas"
$ws.Cells.Item(107,9).Value = 5.25
$ws.Range("H107").WrapText = $true
$ws.Range("I107").WrapText = $true

$ws.Cells.Item(108,1).Value = "DHTD-84"
$ws.Cells.Item(108,2).Value = "Implementation"
$ws.Cells.Item(108,3).Value = "Construction"
$ws.Cells.Item(108,4).Value = "Implement new enemy car"
$ws.Cells.Item(108,6).Value = 4
$ws.Cells.Item(108,7).Value = 2
$ws.Cells.Item(108,8).Value = "Luca"
$ws.Cells.Item(108,9).Value = 0.75
$ws.Range("D108").WrapText = $true
$ws.Range("H108").WrapText = $true
$ws.Range("I108").WrapText = $true

$ws.Cells.Item(109,8).Value = "Nicolas"
$ws.Cells.Item(109,9).Value = 1.25
$ws.Range("H109").WrapText = $true
$ws.Range("I109").WrapText = $true

$ws.Cells.Item(110,1).Value = "DHTD-85"
$ws.Cells.Item(110,2).Value = "Implementation"
$ws.Cells.Item(110,3).Value = "Construction"
$ws.Cells.Item(110,4).Value = "implement new enemy plane"
$ws.Cells.Item(110,6).Value = 3
$ws.Cells.Item(110,7).Value = 2.75
$ws.Cells.Item(110,8).Value = "Fabian"
$ws.Cells.Item(110,9).Value = 2.25
$ws.Range("D110").WrapText = $true
$ws.Range("H110").WrapText = $true
$ws.Range("I110").WrapText = $true

$ws.Cells.Item(111,8).Value = "Luca"
$ws.Cells.Item(111,9).Value = 0.25
$ws.Range("H111").WrapText = $true
$ws.Range("I111").WrapText = $true

$ws.Cells.Item(112,8).Value = "Nicolas"
$ws.Cells.Item(112,9).Value = 0.25
$ws.Range("H112").WrapText = $true
$ws.Range("I112").WrapText = $true

$ws.Cells.Item(113,1).Value = "DHTD-88"
$ws.Cells.Item(113,2).Value = "Implementation"
$ws.Cells.Item(113,3).Value = "Construction"
$ws.Cells.Item(113,4).Value = "Design the easy match"
$ws.Cells.Item(113,6).Value = 4
$ws.Cells.Item(113,7).Value = 2.75
$ws.Cells.Item(113,8).Value = "Fabian"
$ws.Cells.Item(113,9).Value = 0.25
$ws.Range("D113").WrapText = $true
$ws.Range("H113").WrapText = $true
$ws.Range("I113").WrapText = $true

$ws.Cells.Item(114,8).Value = "Luca"
$ws.Cells.Item(114,9).Value = 2.25
$ws.Range("H114").WrapText = $true
$ws.Range("I114").WrapText = $true

$ws.Cells.Item(115,8).Value = "Nicolas"
$ws.Cells.Item(115,9).Value = 0.25
$ws.Range("H115").WrapText = $true
$ws.Range("I115").WrapText = $true

$ws.Cells.Item(116,1).Value = "DHTD-91"
$ws.Cells.Item(116,2).Value = "Implementation"
$ws.Cells.Item(116,3).Value = "Construction"
$ws.Cells.Item(116,4).Value = "Implement boss enemy"
$ws.Cells.Item(116,6).Value = 4
$ws.Cells.Item(116,7).Value = 3.75
$ws.Cells.Item(116,8).Value = "Fabian"
$ws.Cells.Item(116,9).Value = 3.25
$ws.Range("D116").WrapText = $true
$ws.Range("H116").WrapText = $true
$ws.Range("I116").WrapText = $true

$ws.Cells.Item(117,8).Value = "Luca"
$ws.Cells.Item(117,9).Value = 0.25
$ws.Range("H117").WrapText = $true
$ws.Range("I117").WrapText = $true

$ws.Cells.Item(118,8).Value = "Nicolas"
$ws.Cells.Item(118,9).Value = 0.25
$ws.Range("H118").WrapText = $true
$ws.Range("I118").WrapText = $true

$ws.Cells.Item(119,1).Value = "DHTD-106"
$ws.Cells.Item(119,2).Value = "Implementation"
$ws.Cells.Item(119,3).Value = "Construction"
$ws.Cells.Item(119,4).Value = "Bug: Tower Radius is below newly spawned enemies"
$ws.Cells.Item(119,6).Value = 3
$ws.Cells.Item(119,7).Value = 0.25
$ws.Cells.Item(119,8).Value = "Nicolas"
$ws.Cells.Item(119,9).Value = 0.25
$ws.Range("H119").WrapText = $true
$ws.Range("I119").WrapText = $true

$ws.Cells.Item(120,1).Value = "DHTD-114"
$ws.Cells.Item(120,2).Value = "Implementation"
$ws.Cells.Item(120,3).Value = "Construction"
$ws.Cells.Item(120,4).Value = "Ingame Music"
$ws.Cells.Item(120,6).Value = 14
$ws.Cells.Item(120,7).Value = 14
$ws.Cells.Item(120,8).Value = "Fabian"
$ws.Cells.Item(120,9).Value = 14
$ws.Range("D120").WrapText = $true
$ws.Range("H120").WrapText = $true
$ws.Range("I120").WrapText = $true

$ws.Cells.Item(121,1).Value = "DHTD-115"
$ws.Cells.Item(121,2).Value = "Implementation"
$ws.Cells.Item(121,3).Value = "Transition"
$ws.Cells.Item(121,4).Value = "Bugfixes"
$ws.Cells.Item(121,6).Value = 4
$ws.Cells.Item(121,7).Value = 3
$ws.Cells.Item(121,8).Value = "Fabian"
$ws.Cells.Item(121,9).Value = 2
$ws.Range("D121").WrapText = $true
$ws.Range("H121").WrapText = $true
$ws.Range("I121").WrapText = $true

$ws.Cells.Item(122,8).Value = "Luca"
$ws.Cells.Item(122,9).Value = 0.25
$ws.Range("H122").WrapText = $true
$ws.Range("I122").WrapText = $true

$ws.Cells.Item(123,8).Value = "Nicolas"
$ws.Cells.Item(123,9).Value = 0.75
$ws.Range("H123").WrapText = $true
$ws.Range("I123").WrapText = $true

$ws.Cells.Item(124,1).Value = "DHTD-116"
$ws.Cells.Item(124,2).Value = "Project Management"
$ws.Cells.Item(124,3).Value = "Construction"
$ws.Cells.Item(124,4).Value = "Blog Entry Week 18"
$ws.Cells.Item(124,6).Value = 1
$ws.Cells.Item(124,7).Value = 1.5
$ws.Cells.Item(124,8).Value = "Fabian"
$ws.Cells.Item(124,9).Value = 0.25
$ws.Range("D124").WrapText = $true
$ws.Range("H124").WrapText = $true
$ws.Range("I124").WrapText = $true

$ws.Cells.Item(125,8).Value = "Luca"
$ws.Cells.Item(125,9).Value = 0.25
$ws.Range("H125").WrapText = $true
$ws.Range("I125").WrapText = $true

$ws.Cells.Item(126,8).Value = "Nicolas"
$ws.Cells.Item(126,9).Value = 1
$ws.Range("H126").WrapText = $true
$ws.Range("I126").WrapText = $true

$ws.Cells.Item(127,1).Value = "DHTD-117"
$ws.Cells.Item(127,2).Value = "Testing"
$ws.Cells.Item(127,3).Value = "Transition"
$ws.Cells.Item(127,4).Value = "Create question formular for UX test"
$ws.Cells.Item(127,6).Value = 3
$ws.Cells.Item(127,7).Value = 2.25
$ws.Cells.Item(127,8).Value = "Fabian"
$ws.Cells.Item(127,9).Value = 0.25
$ws.Range("D127").WrapText = $true
$ws.Range("H127").WrapText = $true
$ws.Range("I127").WrapText = $true

$ws.Cells.Item(128,8).Value = "Luca"
$ws.Cells.Item(128,9).Value = 1
$ws.Range("H128").WrapText = $true
$ws.Range("I128").WrapText = $true

$ws.Cells.Item(129,8).Value = "Nicolas"
$ws.Cells.Item(129,9).Value = 1
$ws.Range("H129").WrapText = $true
$ws.Range("I129").WrapText = $true

$ws.Cells.Item(130,1).Value = "DHTD-118"
$ws.Cells.Item(130,2).Value = "Analysis & Design"
$ws.Cells.Item(130,3).Value = "Elaboration"
$ws.Cells.Item(130,4).Value = "Create design for plane enemy"
$ws.Cells.Item(130,6).Value = 3
$ws.Cells.Item(130,7).Value = 2.75
$ws.Cells.Item(130,8).Value = "Nicolas"
$ws.Cells.Item(130,9).Value = 2.75
$ws.Range("D130").WrapText = $true
$ws.Range("H130").WrapText = $true
$ws.Range("I130").WrapText = $true

$ws.Range("M17").Formula = "=SUM(G2,G55/2,G100)"
